# aggiornamento fino a 20/09/2021
# Append 11 new daily rows (375-385) to Sheet1, covering 2021-09-10 .. 2021-09-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(44449, 2, 9, 82.91110087517274),
    @(44450, 1, 9, 82.91110087517274),
    @(44451, 4, 11, 101.3357899585444),
    @(44452, 1, 12, 110.5481345002303),
    @(44453, 1, 11, 101.3357899585444),
    @(44454, 0, 10, 92.12344541685859),
    @(44455, 2, 11, 101.3357899585444),
    @(44456, 1, 10, 92.12344541685859),
    @(44457, 0, 9, 82.91110087517274),
    @(44458, 0, 5, 46.06172270842929),
    @(44459, 4, 8, 73.69875633348687)
)

$startRow = 375
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A carries the date style (s="2", thin border, bold, centered, date
# number format) used throughout the sheet - replicate it onto the new cells
# the same way Excel's fill-down / copy-paste-format would.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)
